# "Quitando los campos departamentos innecesarios"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modificacion pre-resolucion")
$ws.Activate()

# Hallazgo #1 ("Implementación de matriz") is now resolved -> mark as closed
$ws.Range("G3").Value = "cerrada"

# Add three more affected tables (views) to the "Tablas Afec" mini-table,
# without the now-unneeded department fields (columns D/E stay untouched)
$ws.Range("H6").Value = "vw_datosrl"
$ws.Range("H7").Value = "vw_accionistas"
$ws.Range("H8").Value = "vw_benefinales"

$ws.Range("I6").Value = "SI"
$ws.Range("J6").Value = "NO"
$ws.Range("K6").Value = "NO"
$ws.Range("L6").Value = "SI"
$ws.Range("M6").Value = "SI"

$ws.Range("I7").Value = "SI"
$ws.Range("J7").Value = "NO"
$ws.Range("K7").Value = "NO"
$ws.Range("L7").Value = "SI"
$ws.Range("M7").Value = "SI"

$ws.Range("I8").Value = "SI"
$ws.Range("J8").Value = "NO"
$ws.Range("K8").Value = "NO"
$ws.Range("L8").Value = "SI"
$ws.Range("M8").Value = "SI"

# Extend the row-grouping merges to cover the new rows
$ws.Range("G3:G8").Merge()
$ws.Range("F3:F8").Merge()
$ws.Range("C3:C8").Merge()
$ws.Range("B3:B8").Merge()
$ws.Range("A3:A8").Merge()

# Normalize row heights on the body rows
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15

# Remove stray underline formatting left on K4
$ws.Range("K4").Font.Underline = 0

# Update selection to reflect where the user ended up
$ws.Range("H9").Select()

$ws2 = $wb.Worksheets.Item("hallazgos")
$ws2.Activate()
$ws2.Range("F5").Select()
